# Fruta / hortaliza, semanal
# Insert a new weekly record as row 7, pushing the existing rows 7:100 down
# to rows 8:101 (the sheet keeps its historical rows, it just gets a new
# entry at the top of the data block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7:100 down to 8:101, carrying formatting/styles with them.
$ws.Rows.Item(7).EntireRow.Insert()

# Populate the newly freed row 7 with the new weekly observation.
$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "Vega Modelo de Temuco"
$ws.Range("C7").Value = "La Araucanía"
$ws.Range("D7").Value = 44881
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100101
$ws.Range("H7").Value = "Berries"
$ws.Range("I7").Value = 100101001
$ws.Range("J7").Value = "Arándano (blue)"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 400
$ws.Range("N7").Value = 3200
$ws.Range("O7").Value = 3200
$ws.Range("P7").Value = 3200
$ws.Range("Q7").Value = "$/kilo"
$ws.Range("R7").Value = "Región del Maule"
$ws.Range("S7").Value = 3200
$ws.Range("T7").Value = 1
